$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.209.41'
$ws.Range('E2').Value = '  +0.47%  '
$ws.Range('D3').Value = '3.145.92'
$ws.Range('E3').Value = '  +0.84%  '
$ws.Range('E4').Value = '  +0.11%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '591.95'
$ws.Range('D5').Style = "Normal"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '145.84'
$ws.Range('D6').Style = "Normal"
$ws.Range('E6').Value = '  -0.26%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.998'
$ws.Range('D7').Style = "Normal"
$ws.Range('E7').Value = '  -0.19%  '
$ws.Range('D8').Value = '3.136.27'
$ws.Range('E8').Value = '  +0.69%  '
$ws.Range('E9').Value = '  -0.06%  '
$ws.Range('E10').Value = '  +0.63%  '
$ws.Range('E11').Value = '  +3.65%  '
$ws.Range('E12').Value = '  -1.12%  '
$ws.Range('E13').Value = '  -1.28%  '
$ws.Range('E14').Value = '  +0.11%  '
$ws.Range('D15').Value = '3.663.57'
$ws.Range('E15').Value = '  +0.75%  '
$ws.Range('E16').Value = '  -0.98%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '7.31'
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = '  +2.53%  '
$ws.Range('D18').Value = '63.990.77'
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('D19').Value = '3.139.33'
$ws.Range('E19').Value = '  +0.60%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '469.80'
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = '  +1.38%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '14.38'
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = '  +0.44%  '
$ws.Range('E22').Value = '  +0.69%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.59'
$ws.Range('D23').Style = "Normal"
$ws.Range('E23').Value = '  +0.85%  '
$ws.Range('E24').Value = '  +9.16%  '
$ws.Range('E25').Value = '  -1.30%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '81.45'
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = '  -0.32%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '0.999'
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = '  -0.07%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '9.86'
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = '  +10.27%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '7.49'
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = '  +8.86%  '
$ws.Range('E30').Value = '  +0.98%  '
$ws.Range('E31').Value = '  +0.57%  '
$ws.Range('E32').Value = '  +0.18%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '27.70'
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = '  +2.74%  '
$ws.Range('E34').Value = '  +1.18%  '
$ws.Range('D35').Value = '0.0₃0851'
$ws.Range('E35').Value = '  -1.91%  '
$ws.Range('E36').Value = '  +1.51%  '
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '6.16'
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = '  +1.73%  '
$ws.Range('E38').Value = '  -2.34%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '3.22'
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = '  -5.15%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '51.38'
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = '  +0.76%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '9.34'
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = '  +7.40%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '454.50'
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = '  +1.86%  '
$ws.Range('E43').Value = '  +6.67%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '0.0374'
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = '  +0.65%  '
$ws.Range('D45').Value = '2.917.14'
$ws.Range('E45').Value = '  +1.70%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '40.62'
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = '  +13.42%  '
$ws.Range('E47').Value = '  -2.52%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '134.36'
$ws.Range('D48').Style = "Normal"
$ws.Range('E48').Value = '  +9.01%  '
$ws.Range('E49').Value = '  +0.00%  '
$ws.Range('E51').Value = '  +3.11%  '
